# Removed investor creation in PI and valuations import
#
# The Portfolio Investments import template used to also collect the
# investor's "Pan" and "Primary Email *" (columns C and D), which the
# sheet used to auto-create an investor record (hence the mailto:
# hyperlink on the sample email in D2). That behaviour was removed, so
# the template columns are no longer needed.
#
# Remove the hyperlink on the sample e-mail first (so the workbook
# doesn't keep a relationship pointing at a cell that is about to be
# deleted), then delete the entire "Pan" / "Primary Email *" columns
# (C:D), shifting every column to their right one left by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Delete()
$ws.Range("C:D").EntireColumn.Delete()

$ws.Range("C1:D1048576").Select()
